$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-14 Sunday", "2025-09-15 Monday"),
    @("428÷3=", "121÷3="),
    @("956÷5=", "691÷5="),
    @("314÷7=", "856÷4="),
    @("922÷5=", "522÷7="),
    @("271÷9=", "170÷7="),
    @("787÷8=", "961÷2="),
    @("564÷6=", "188÷4="),
    @("153÷9=", "229÷4="),
    @("576÷3=", "965÷6="),
    @("903÷6=", "583÷6="),
    @("480÷7=", "601÷8="),
    @("630÷7=", "345÷4="),
    @("733÷5=", "127÷4="),
    @("118÷8=", "886÷9="),
    @("923÷7=", "422÷8="),
    @("213÷7=", "760÷2="),
    @("366÷6=", "815÷3="),
    @("416÷9=", "249÷6="),
    @("858÷6=", "538÷5="),
    @("750÷7=", "896÷8="),
    @("669÷9=", "556÷5="),
    @("501÷7=", "127÷2="),
    @("856÷9=", "189÷4="),
    @("109÷8=", "558÷5="),
    @("129÷2=", "622÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
